$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for existing data rows 2-37
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45696
}

# Set explicit height on row 37 (matches Excel behavior when re-saving)
$ws.Rows.Item(37).RowHeight = 15

# Append new row 38 with the new entry
$ws.Cells.Item(38, 1).Value2 = "A 6050-2025"
$ws.Cells.Item(38, 2).Value2 = 45695
$ws.Cells.Item(38, 3).Value2 = 45696
$ws.Cells.Item(38, 4).Value2 = "OKÄNT"
$ws.Cells.Item(38, 5).Value2 = "OKÄNT"
$ws.Cells.Item(38, 7).Value2 = 5.5
$ws.Cells.Item(38, 8).Value2 = 0
$ws.Cells.Item(38, 9).Value2 = 0
$ws.Cells.Item(38, 10).Value2 = 0
$ws.Cells.Item(38, 11).Value2 = 0
$ws.Cells.Item(38, 12).Value2 = 0
$ws.Cells.Item(38, 13).Value2 = 0
$ws.Cells.Item(38, 14).Value2 = 0
$ws.Cells.Item(38, 15).Value2 = 0
$ws.Cells.Item(38, 16).Value2 = 0
$ws.Cells.Item(38, 17).Value2 = 0

# Apply same number format as dates column to B38/C38
$ws.Cells.Item(38, 2).NumberFormat = $ws.Cells.Item(37, 2).NumberFormat
$ws.Cells.Item(38, 3).NumberFormat = $ws.Cells.Item(37, 3).NumberFormat

# R38 should mimic R37's style (wrap text, empty text)
$ws.Cells.Item(38, 18).Value2 = "'"
$ws.Cells.Item(38, 18).WrapText = $true
